# "windowstate en windowstyle transferable"
# Update the "2Player1" highscore sheet: drop the "Tjitske" entry and
# renumber/extend the points column (B4:B6), keeping the existing
# descending sort on column B in sync with the new data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2Player1")

# Row 4 used to hold "Tjitske" / 50. Remove the name and fold the row's
# points into the now-longer points list (B4:B6 = 300, 300, 100).
$ws.Range("A4").ClearContents()
$ws.Range("B4").Value = 300
$ws.Range("B5").Value = 300
$ws.Range("B6").Value = 100

# Keep the sheet's sort definition in sync with the new A3:B6 data range
# (still sorted descending on the Points column), so the dimension/sort
# metadata matches the refreshed data.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B3:B6"), 0, 2) | Out-Null
$ws.Sort.SetRange($ws.Range("A3:B6"))
$ws.Sort.Header = 0
$ws.Sort.Apply()
